$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.163197755813599
$ws.Range("B1").Value = 1.526278138160706
$ws.Range("C1").Value = 1.548879742622375
$ws.Range("D1").Value = 1.470386505126953
$ws.Range("E1").Value = 1.386083006858826
